$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "선형 연산자와 함수 공간"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/05/31/linear_operator_and_function_space.html"

$ws.Range("D9").Value = "MBA in AI BigData만 해도 업계 상위 1%일텐데요?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/mba-in-ai-bigdata-top-1-percent/#utm_source=rss&utm_medium=rss&utm_campaign=mba-in-ai-bigdata-top-1-percent"

$ws.Range("D51").Value = "[MariaDB] Group by 로 그룹화된 그룹의 갯수"
$ws.Range("E51").Value = "https://bskyvision.com/1195"
